# Weekly price-update edit: a new daily price record for "Cebollín" at
# "Vega Modelo de Temuco" is inserted as row 174 (pushing the existing
# rows 174-252 down to 175-253).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 174; everything below shifts down one row.
$ws.Rows("174").Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A174").Value = 10
$ws.Range("B174").Value = "Vega Modelo de Temuco"
$ws.Range("C174").Value = "La Araucanía"
$ws.Range("D174").Value = 44518
$ws.Range("E174").Value = 9
$ws.Range("F174").Value = 100112037
$ws.Range("G174").Value = "Cebollín"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 50
$ws.Range("K174").Value = 8000
$ws.Range("L174").Value = 8000
$ws.Range("M174").Value = 8000
$ws.Range("N174").Value = "`$/docena de paquetes"
$ws.Range("O174").Value = "Provincia de Cautín"
$ws.Range("P174").Value = 667
$ws.Range("Q174").Value = 12
$ws.Range("R174").Value = "Hortaliza"
